$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) to short machine-readable names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish connector words (de/del/la/las/los/el/y) in state/municipality names ---
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B27").Value = "Amatenango De La Frontera"
$ws.Range("B43").Value = "San Cristóbal De Las Casas"
$ws.Range("B79").Value = "Guadalupe Y Calvo"
$ws.Range("B82").Value = "Hidalgo Del Parral"
$ws.Range("B104").Value = "San Francisco De Borja"
$ws.Range("B105").Value = "San Francisco De Conchos"
$ws.Range("B106").Value = "San Francisco Del Oro"
$ws.Range("B114").Value = "Valle De Zaragoza"
$ws.Range("B125").Value = "San Juan De Sabinas"
$ws.Range("A136").Value = "Ciudad De México"
$ws.Range("B140").Value = "Cuajimalpa De Morelos"
$ws.Range("B154").Value = "Coneto De Comonfort"
$ws.Range("B168").Value = "Nombre De Dios"
$ws.Range("B172").Value = "Pánuco De Coronado"
$ws.Range("B179").Value = "San Juan De Guadalupe"
$ws.Range("B180").Value = "San Juan Del Río"
$ws.Range("B181").Value = "San Luis Del Cordero"
$ws.Range("B182").Value = "San Pedro Del Gallo"
$ws.Range("A191").Value = "Estado De México"
$ws.Range("B191").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B197").Value = "Coacalco De Berriozábal"
$ws.Range("B203").Value = "Ecatepec De Morelos"
$ws.Range("B211").Value = "Naucalpan De Juárez"
$ws.Range("B221").Value = "Tlalnepantla De Baz"
$ws.Range("B223").Value = "Valle De Bravo"
$ws.Range("B231").Value = "San Miguel De Allende"
$ws.Range("B232").Value = "Apaseo El Alto"
$ws.Range("B239").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B243").Value = "Jaral Del Progreso"
$ws.Range("B250").Value = "Purísima Del Rincón"
$ws.Range("B254").Value = "San Diego De La Unión"
$ws.Range("B256").Value = "San Francisco Del Rincón"
$ws.Range("B258").Value = "San Luis De La Paz"
$ws.Range("B259").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B260").Value = "Silao De La Victoria"
$ws.Range("B264").Value = "Valle De Santiago"
$ws.Range("B269").Value = "Acapulco De Juárez"
$ws.Range("B270").Value = "Ajuchitlán Del Progreso"
$ws.Range("B271").Value = "Atenango Del Río"
$ws.Range("B272").Value = "Atoyac De Álvarez"
$ws.Range("B273").Value = "Ayutla De Los Libres"
$ws.Range("B275").Value = "Chilpancingo De Los Bravo"
$ws.Range("B277").Value = "Coyuca De Benítez"
$ws.Range("B278").Value = "Coyuca De Catalán"
$ws.Range("B279").Value = "Cuetzala Del Progreso"
$ws.Range("B280").Value = "Cutzamala De Pinzón"
$ws.Range("B284").Value = "Huitzuco De Los Figueroa"
$ws.Range("B285").Value = "Iguala De La Independencia"
$ws.Range("B286").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B287").Value = "Zihuatanejo De Azueta"
$ws.Range("B294").Value = "Taxco De Alarcón"
$ws.Range("B296").Value = "Técpan De Galeana"
$ws.Range("B304").Value = "Atotonilco El Grande"
$ws.Range("B306").Value = "Cuautepec De Hinojosa"
$ws.Range("B311").Value = "Mixquiahuala De Juárez"
$ws.Range("B312").Value = "Pachuca De Soto"
$ws.Range("B316").Value = "Tepehuacán De Guerrero"
$ws.Range("B317").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B319").Value = "Tula De Allende"
$ws.Range("B320").Value = "Tulancingo De Bravo"
$ws.Range("B324").Value = "Zapotlán De Juárez"
$ws.Range("B328").Value = "Acatlán De Juárez"
$ws.Range("B329").Value = "Ahualulco De Mercado"
$ws.Range("B333").Value = "Atotonilco El Alto"
$ws.Range("B334").Value = "Autlán De Navarro"
$ws.Range("B345").Value = "Encarnación De Díaz"
$ws.Range("B348").Value = "Huejuquilla El Alto"
$ws.Range("B349").Value = "Ixtlahuacán Del Río"
$ws.Range("B352").Value = "Lagos De Moreno"
$ws.Range("B357").Value = "Ojuelos De Jalisco"
$ws.Range("B361").Value = "San Cristóbal De La Barranca"
$ws.Range("B362").Value = "San Diego De Alejandría"
$ws.Range("B364").Value = "San Juan De Los Lagos"
$ws.Range("B367").Value = "San Martín De Bolaños"
$ws.Range("B368").Value = "San Miguel El Alto"
$ws.Range("B369").Value = "Santa María De Los Ángeles"
$ws.Range("B372").Value = "Tamazula De Gordiano"
$ws.Range("B377").Value = "Teocuitatlán De Corona"
$ws.Range("B378").Value = "Tepatitlán De Morelos"
$ws.Range("B380").Value = "Tizapán El Alto"
$ws.Range("B381").Value = "Tlajomulco De Zúñiga"
$ws.Range("B386").Value = "Unión De Tula"
$ws.Range("B391").Value = "Yahualica De González Gallo"
$ws.Range("B394").Value = "Zapotlán Del Rey"
$ws.Range("B395").Value = "Zapotlán El Grande"
$ws.Range("B459").Value = "Tlaltizapán De Zapata"
$ws.Range("B466").Value = "Ixtlán Del Río"
$ws.Range("B471").Value = "Santa María Del Oro"
$ws.Range("B483").Value = "San Nicolás De Los Garza"
$ws.Range("B488").Value = "Guevea De Humboldt"
$ws.Range("B489").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B490").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B491").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B492").Value = "Ixtlán De Juárez"
$ws.Range("B493").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B495").Value = "Mariscala De Juárez"
$ws.Range("B497").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B498").Value = "Oaxaca De Juárez"
$ws.Range("B499").Value = "Ocotlán De Morelos"
$ws.Range("B517").Value = "Santo Domingo De Morelos"
$ws.Range("B520").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B521").Value = "Tlacolula De Matamoros"
$ws.Range("B522").Value = "Villa De Chilapa De Díaz"
$ws.Range("B523").Value = "Villa Sola De Vega"
$ws.Range("B529").Value = "Chalchicomula De Sesma"
$ws.Range("B535").Value = "Cuayuca De Andrade"
$ws.Range("B541").Value = "Izúcar De Matamoros"
$ws.Range("B546").Value = "Palmar De Bravo"
$ws.Range("B556").Value = "Tepexi De Rodríguez"
$ws.Range("B557").Value = "Tetela De Ocampo"
$ws.Range("B566").Value = "Cadereyta De Montes"
$ws.Range("B569").Value = "Jalpan De Serra"
$ws.Range("B573").Value = "San Juan Del Río"
$ws.Range("B581").Value = "Ciudad Del Maíz"
$ws.Range("B586").Value = "Mexquitic De Carmona"
$ws.Range("B593").Value = "Santa María Del Río"
$ws.Range("B595").Value = "Soledad De Graciano Sánchez"
$ws.Range("B598").Value = "Villa De Ramos"
$ws.Range("B633").Value = "Nacozari De García"
$ws.Range("B661").Value = "Contla De Juan Cuamatzi"
$ws.Range("B663").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B679").Value = "Cosamaloapan De Carpio"
$ws.Range("B689").Value = "Martínez De La Torre"
$ws.Range("B690").Value = "Medellín De Bravo"
$ws.Range("B699").Value = "Paso De Ovejas"
$ws.Range("B702").Value = "Poza Rica De Hidalgo"
$ws.Range("B718").Value = "Cañitas De Felipe Pescador"
$ws.Range("B720").Value = "Concepción Del Oro"
$ws.Range("B735").Value = "Moyahua De Estrada"
$ws.Range("B736").Value = "Nochistlán De Mejía"
$ws.Range("B737").Value = "Noria De Ángeles"
$ws.Range("B746").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B748").Value = "Villa De Cos"

# --- Remove trailing metadata/footnote rows (756-760): sample size, source, author, date ---
$ws.Range("A756:D760").EntireRow.Delete()

